$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '42.702.95'
$ws.Cells.Item(2, 5).Value = '  +1.40%  '
$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '2.304.36'
$ws.Cells.Item(3, 5).Value = '  +0.76%  '
$ws.Cells.Item(4, 5).Value = '  +0.15%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '316.86'
$ws.Cells.Item(5, 5).Value = '  +0.16%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '103.81'
$ws.Cells.Item(6, 5).Value = '  +0.53%  '
$ws.Cells.Item(8, 5).Value = '  +0.29%  '
$ws.Cells.Item(9, 5).Value = '  +0.70%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '39.96'
$ws.Cells.Item(10, 5).Value = '  +1.84%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0905'
$ws.Cells.Item(11, 5).Value = '  +0.34%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '8.51'
$ws.Cells.Item(12, 5).Value = '  +3.41%  '
$ws.Cells.Item(13, 5).Value = '  +1.06%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.995'
$ws.Cells.Item(14, 5).Value = '  +3.97%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '15.37'
$ws.Cells.Item(15, 5).Value = '  +1.33%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '2.654.65'
$ws.Cells.Item(16, 5).Value = '  +0.84%  '
$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '2.302.90'
$ws.Cells.Item(17, 5).Value = '  +0.75%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '42.636.82'
$ws.Cells.Item(18, 5).Value = '  +1.53%  '
$ws.Cells.Item(19, 5).Value = '  +3.43%  '
$ws.Cells.Item(20, 5).Value = '  +0.86%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '13.49'
$ws.Cells.Item(21, 5).Value = '  +33.82%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '74.07'
$ws.Cells.Item(22, 5).Value = '  +1.31%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '3.55'
$ws.Cells.Item(23, 5).Value = '  -1.81%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '268.03'
$ws.Cells.Item(24, 5).Value = '  -3.67%  '
$ws.Cells.Item(25, 5).Value = '  -0.50%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '1.01'
$ws.Cells.Item(26, 5).Value = '  -0.18%  '
$ws.Cells.Item(27, 5).Value = '  +1.65%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '2.34'
$ws.Cells.Item(28, 5).Value = '  -2.78%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '22.63'
$ws.Cells.Item(29, 5).Value = '  -0.63%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '38.10'
$ws.Cells.Item(30, 5).Value = '  +6.21%  '
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '6.53'
$ws.Cells.Item(31, 5).Value = '  +12.47%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '166.09'
$ws.Cells.Item(32, 5).Value = '  +2.05%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.0882'
$ws.Cells.Item(33, 5).Value = '  +1.74%  '
$ws.Cells.Item(34, 2).Value = 'Stellar'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.132'
$ws.Cells.Item(34, 5).Value = '  -3.29%  '
$ws.Cells.Item(35, 2).Value = 'WEMIXToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '2.66'
$ws.Cells.Item(35, 5).Value = '  -6.53%  '
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.114'
$ws.Cells.Item(36, 5).Value = '  +0.29%  '
$ws.Cells.Item(37, 5).Value = '  +2.24%  '
$ws.Cells.Item(38, 5).Value = '  +2.36%  '
$ws.Cells.Item(39, 2).Value = 'LidoDAOToken'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '2.78'
$ws.Cells.Item(39, 5).Value = '  -1.98%  '
$ws.Cells.Item(40, 2).Value = 'NEARProtocol'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '3.70'
$ws.Cells.Item(40, 5).Value = '  -1.05%  '
$ws.Cells.Item(41, 5).Value = '  +13.99%  '
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '97.69'
$ws.Cells.Item(42, 5).Value = '  -1.61%  '
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '69.96'
$ws.Cells.Item(43, 5).Value = '  +1.38%  '
$ws.Cells.Item(44, 5).Value = '  +1.04%  '
$ws.Cells.Item(45, 5).Value = '  -0.07%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '118.17'
$ws.Cells.Item(46, 5).Value = '  +5.18%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '12.34'
$ws.Cells.Item(47, 5).Value = '  +3.99%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '79.99'
$ws.Cells.Item(48, 5).Value = '  +4.02%  '
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '1.650.38'
$ws.Cells.Item(49, 5).Value = '  +4.69%  '
$ws.Cells.Item(50, 5).Value = '  +0.81%  '
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '8.87'
$ws.Cells.Item(51, 5).Value = '  -0.03%  '
